$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "non_tumor" was renamed to "normal" throughout the Tissue_type column.
$ws.Cells.Replace("non_tumor", "normal")

# Fill in the previously-blank Tissue_type values for rows 72-81 with
# "chronic_pancreatitis", matching the surrounding rows for that study.
for ($r = 72; $r -le 81; $r++) {
    $ws.Cells.Item($r, 2).Value = "chronic_pancreatitis"
}
